$d = $word.ActiveDocument

# 1. Remove the stale "_GoBack" bookmark (it currently sits in the empty
#    paragraph right after the introduction text). Word will recreate this
#    bookmark automatically at the location of the next real edit.
$d.Bookmarks("_GoBack").Delete()

# 2. Merge the two list paragraphs
#      "De extra database voorbereidingen dienen niet gebruikt te worden in "
#      "de applicatie (code)."
#    back into a single list item by deleting the paragraph break between
#    them.
$r = $d.Content
$r.Find.Execute("De extra database voorbereidingen dienen niet gebruikt te worden in ")
$mergePoint = $r.End
$d.Range($mergePoint, $mergePoint + 1).Delete()

# 3. Re-create "_GoBack" at the spot of this edit (between the two runs that
#    were just joined into one paragraph) - this is exactly what Word does
#    whenever new text is edited/deleted.
$d.Bookmarks.Add("_GoBack", $d.Range($mergePoint, $mergePoint))

# 4. The page count shown in the footer moved from 4 to 8 pages.
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("4", $false, $false, $false, $false, $false, $true, 1, $false, "8", 2)
